$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 587
$ws.Range("I4").Value = 587
$ws.Range("K4").Value = 587
$ws.Range("M4").Value = -473

$ws.Range("H18").Value = 2004.75
$ws.Range("I18").Value = 2004.75
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2004.75
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1720.75
$ws.Range("N18").ClearContents()

$ws.Range("H58").Value = 1618.8889
$ws.Range("I58").Value = 94
$ws.Range("J58").Value = 3525
$ws.Range("K58").Value = 282
$ws.Range("L58").Value = 10575
$ws.Range("M58").Value = -132
$ws.Range("N58").Value = -10875

$ws.Range("H61").Value = 433.16666
$ws.Range("I61").Value = 187.5
$ws.Range("J61").Value = 924.5
$ws.Range("K61").Value = 562.5
$ws.Range("L61").Value = 2773.5
$ws.Range("M61").Value = -390.5
$ws.Range("N61").Value = -3117.5

$ws.Range("H62").Value = 3855.5
$ws.Range("I62").Value = 3142
$ws.Range("K62").Value = 3142
$ws.Range("M62").Value = -2518

$ws.Range("H65").Value = 3855.5
$ws.Range("I65").Value = 3142
$ws.Range("K65").Value = 15710
$ws.Range("M65").Value = -12590

$ws.Range("H74").Value = 9048.909
$ws.Range("I74").Value = 6907.8
$ws.Range("K74").Value = 6907.8
$ws.Range("M74").Value = -5971.8

$ws.Range("H77").Value = 9048.909
$ws.Range("I77").Value = 6907.8
$ws.Range("K77").Value = 34539
$ws.Range("M77").Value = -29859

$ws.Range("H106").Value = 1302.4
$ws.Range("I106").Value = 1373
$ws.Range("K106").Value = 1373
$ws.Range("M106").Value = -742

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H131").Value = 5206.5835
$ws.Range("I131").Value = 2497.6667
$ws.Range("K131").Value = 7493.000100000001
$ws.Range("M131").Value = -2453.000100000001

$ws.Range("H132").Value = 16501.924
$ws.Range("I132").Value = 2618.8667
$ws.Range("K132").Value = 7856.6001
$ws.Range("M132").Value = -5326.6001

$ws.Range("H137").Value = 2497.5745
$ws.Range("I137").Value = 2330.8572
$ws.Range("J137").Value = 3898
$ws.Range("K137").Value = 6992.571599999999
$ws.Range("L137").Value = 11694
$ws.Range("M137").Value = -4442.571599999999
$ws.Range("N137").Value = -16794

$ws.Range("H138").Value = 4197.431
$ws.Range("J138").Value = 4460.451
$ws.Range("L138").Value = 13381.353
$ws.Range("N138").Value = -23661.353

$ws.Range("H141").Value = 3942.889
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1270.94
$ws.Range("I32").Value = 1279.2245
$ws.Range("K32").Value = 1279.2245
$ws.Range("M32").Value = -992.2245

$ws.Range("H61").Value = 2582.3667
$ws.Range("I61").Value = 2454.28
$ws.Range("K61").Value = 2454.28
$ws.Range("M61").Value = -2242.28

$ws.Range("H74").Value = 2711.3333
$ws.Range("I74").Value = 3042.25
$ws.Range("J74").Value = 1869
$ws.Range("K74").Value = 3042.25
$ws.Range("L74").Value = 1869
$ws.Range("M74").Value = -2168.25
$ws.Range("N74").Value = -3617

$ws.Range("H77").Value = 2711.3333
$ws.Range("I77").Value = 3042.25
$ws.Range("J77").Value = 1869
$ws.Range("K77").Value = 15211.25
$ws.Range("L77").Value = 9345
$ws.Range("M77").Value = -10843.25
$ws.Range("N77").Value = -18081

$ws.Range("H132").Value = 2768.6086
$ws.Range("I132").Value = 2378.7334
$ws.Range("K132").Value = 7136.2002
$ws.Range("M132").Value = -4606.2002

$ws.Range("H136").Value = 2582.3667
$ws.Range("I136").Value = 2454.28
$ws.Range("K136").Value = 7362.84
$ws.Range("M136").Value = -4812.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 940.52
$ws.Range("I22").Value = 563.7273
$ws.Range("J22").Value = 1236.5714
$ws.Range("K22").Value = 563.7273
$ws.Range("L22").Value = 1236.5714
$ws.Range("M22").Value = -213.7273
$ws.Range("N22").Value = -1936.5714

$ws.Range("H31").Value = 2151.2222
$ws.Range("I31").Value = 1950.5518
$ws.Range("J31").Value = 2982.5715
$ws.Range("K31").Value = 1950.5518
$ws.Range("L31").Value = 2982.5715
$ws.Range("M31").Value = -1655.5518
$ws.Range("N31").Value = -3572.5715

$ws.Range("H34").Value = 2151.2222
$ws.Range("I34").Value = 1950.5518
$ws.Range("J34").Value = 2982.5715
$ws.Range("K34").Value = 1950.5518
$ws.Range("L34").Value = 2982.5715
$ws.Range("M34").Value = -1748.5518
$ws.Range("N34").Value = -3386.5715

$ws.Range("H58").Value = 2110.7
$ws.Range("I58").Value = 1230.381
$ws.Range("J58").Value = 4164.778
$ws.Range("K58").Value = 1230.381
$ws.Range("L58").Value = 4164.778
$ws.Range("M58").Value = -1027.381
$ws.Range("N58").Value = -4570.778

$ws.Range("H99").Value = 10203050
$ws.Range("I99").Value = 2219058
$ws.Range("J99").Value = 18187042
$ws.Range("K99").Value = 2219058
$ws.Range("L99").Value = 18187042
$ws.Range("M99").Value = -2217560
$ws.Range("N99").Value = -18190038

$ws.Range("H126").Value = 10203050
$ws.Range("I126").Value = 2219058
$ws.Range("J126").Value = 18187042
$ws.Range("K126").Value = 6657174
$ws.Range("L126").Value = 54561126
$ws.Range("M126").Value = -6654704
$ws.Range("N126").Value = -54566066

$ws.Range("H132").Value = 1996.1
$ws.Range("I132").Value = 1996.1
$ws.Range("K132").Value = 5988.299999999999
$ws.Range("M132").Value = -3458.299999999999

$ws.Range("H134").Value = 3299.4075
$ws.Range("J134").Value = 4638.875
$ws.Range("L134").Value = 13916.625
$ws.Range("N134").Value = -18986.625

$ws.Range("H136").Value = 2110.7
$ws.Range("I136").Value = 1230.381
$ws.Range("J136").Value = 4164.778
$ws.Range("K136").Value = 3691.143
$ws.Range("L136").Value = 12494.334
$ws.Range("M136").Value = -1141.143
$ws.Range("N136").Value = -17594.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 657.36365
$ws.Range("I23").Value = 1865
$ws.Range("J23").Value = 204.5
$ws.Range("K23").Value = 5595
$ws.Range("L23").Value = 613.5
$ws.Range("M23").Value = -5360
$ws.Range("N23").Value = -1083.5

$ws.Range("H29").Value = 4575552
$ws.Range("I29").Value = 8642445
$ws.Range("K29").Value = 25927335
$ws.Range("M29").Value = -25927058

$ws.Range("H33").Value = 20.5
$ws.Range("I33").Value = 20.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 123
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 160
$ws.Range("N33").ClearContents()

$ws.Range("H109").Value = 567.9091
$ws.Range("I109").Value = 424.7
$ws.Range("K109").Value = 1274.1
$ws.Range("M109").Value = -234.0999999999999

$ws.Range("H137").Value = 3214.85
$ws.Range("J137").Value = 3483.0833
$ws.Range("L137").Value = 10449.2499
$ws.Range("N137").Value = -20649.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 335.25
$ws.Range("I2").Value = 137.28572
$ws.Range("J2").Value = 489.22223
$ws.Range("K2").Value = 137.28572
$ws.Range("L2").Value = 489.22223
$ws.Range("M2").Value = -24.28572
$ws.Range("N2").Value = -715.2222300000001

$ws.Range("H102").Value = 10624.5
$ws.Range("I102").Value = 6750
$ws.Range("J102").Value = 14499
$ws.Range("K102").Value = 6750
$ws.Range("L102").Value = 14499
$ws.Range("M102").Value = -5128
$ws.Range("N102").Value = -17743

$ws.Range("H132").Value = 1913.3077
$ws.Range("I132").Value = 1739.4166
$ws.Range("K132").Value = 5218.2498
$ws.Range("M132").Value = -2688.2498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 564.24
$ws.Range("I16").Value = 695.6667
$ws.Range("J16").Value = 226.28572
$ws.Range("K16").Value = 695.6667
$ws.Range("L16").Value = 226.28572
$ws.Range("M16").Value = -525.6667
$ws.Range("N16").Value = -566.28572

$ws.Range("H22").Value = 2149.6667
$ws.Range("J22").Value = 2179.6
$ws.Range("L22").Value = 2179.6
$ws.Range("N22").Value = -2769.6

$ws.Range("H27").Value = 2149.6667
$ws.Range("J27").Value = 2179.6
$ws.Range("L27").Value = 2179.6
$ws.Range("N27").Value = -2393.6

$ws.Range("H82").Value = 1991.2778
$ws.Range("I82").Value = 2913.8572
$ws.Range("J82").Value = 1404.1818
$ws.Range("K82").Value = 2913.8572
$ws.Range("L82").Value = 1404.1818
$ws.Range("M82").Value = -2552.8572
$ws.Range("N82").Value = -2126.1818

$ws.Range("H85").Value = 1991.2778
$ws.Range("I85").Value = 2913.8572
$ws.Range("J85").Value = 1404.1818
$ws.Range("K85").Value = 2913.8572
$ws.Range("L85").Value = 1404.1818
$ws.Range("M85").Value = -1665.8572
$ws.Range("N85").Value = -3900.1818

$ws.Range("H132").Value = 3783.3948
$ws.Range("I132").Value = 2452
$ws.Range("J132").Value = 8776.125
$ws.Range("K132").Value = 7356
$ws.Range("L132").Value = 26328.375
$ws.Range("M132").Value = -4826
$ws.Range("N132").Value = -31388.375

$ws.Range("H136").Value = 1932.0889
$ws.Range("I136").Value = 1762.7949
$ws.Range("K136").Value = 5288.384700000001
$ws.Range("M136").Value = -2738.384700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1547.3334
$ws.Range("J126").Value = 1660.5
$ws.Range("L126").Value = 4981.5
$ws.Range("N126").Value = -9921.5

$ws.Range("H136").Value = 1889.7693
$ws.Range("J136").Value = 4354.9
$ws.Range("L136").Value = 13064.7
$ws.Range("M136").Value = -742
$ws.Range("N136").Value = -18164.7

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
